$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (DAMSLTag, DialogAct)
$updates = @{
    9   = @("sd", "Statement-non-opinion")
    13  = @("sd", "Statement-non-opinion")
    22  = @("sd", "Statement-non-opinion")
    23  = @("aa", "Agree/Accept")
    24  = @("sv", "Statement-opinion")
    32  = @("%", "Uninterpretable")
    37  = @("b", "Acknowledge (Backchannel)")
    42  = @("b", "Acknowledge (Backchannel)")
    55  = @("sd", "Statement-non-opinion")
    95  = @("sv", "Statement-opinion")
    119 = @("aa", "Agree/Accept")
    120 = @("aa", "Agree/Accept")
    127 = @("sd", "Statement-non-opinion")
    134 = @("sd", "Statement-non-opinion")
    144 = @("%", "Uninterpretable")
    169 = @("ba", "Appreciation")
    173 = @("sd", "Statement-non-opinion")
    211 = @("b", "Acknowledge (Backchannel)")
    224 = @("b", "Acknowledge (Backchannel)")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
